# Update "南宁-漫展信息.xlsx" to match output generated at 456a3b4:
#  - sheet "展览" (index 1) and sheet "全部类型" (index 4) both get:
#      * F3: 8376 -> 8398 (想去人数 for the 2024-02-15 event)
#      * F5: 354  -> 358  (想去人数 for the 2024-03-16 event)
#      * a new row 6 appended with the 2024-03-30 event

$wb = $excel.ActiveWorkbook

foreach ($idx in 1, 4) {
    $ws = $wb.Worksheets.Item($idx)

    # Updated "want to go" counts for existing rows.
    $ws.Range("F3").Value = 8398
    $ws.Range("F5").Value = 358

    # New row 6: copy the formatting of row 5's index cell (bold/border/
    # centered style) onto the new index cell, then fill in the values.
    $ws.Range("A5").Copy($ws.Range("A6"))

    $ws.Range("A6").Value = 5

    # "2024-03-30" looks like a date literal, so Excel's COM layer would
    # silently convert it to a date serial unless it is entered as text
    # (leading apostrophe = force-text entry, same as typing it by hand);
    # resetting the style to "Normal" afterwards drops the quote-prefix
    # formatting flag so the cell ends up with the plain default style,
    # matching the other text cells in the sheet.
    $ws.Range("B6").Value = "'2024-03-30"
    $ws.Range("B6").Style = "Normal"

    $ws.Range("C6").Value = "南宁·第一届ANE·DACG动漫嘉年华"
    $ws.Range("D6").Value = "亭洪路45号 百益上河城"
    $ws.Range("E6").Value = "2024.03.30 09:00-03.31 17:30"
    $ws.Range("F6").Value = 0
    $ws.Range("G6").Value = 60
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81658"
    $ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202402/3syliqwc1706852024845.jpeg"
}
